# Applies the cryptos list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.002.62"
$ws.Range("E2").Value = "  +1.53%  "

# Row 3
$ws.Range("D3").Value = "3.131.74"
$ws.Range("E3").Value = "  +1.14%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "599.63"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "141.90"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.55%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "3.120.77"
$ws.Range("E8").Value = "  +0.88%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.519"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.85%  "

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.148"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.05%  "

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.32"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +2.64%  "

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.466"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.39%  "

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000251"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.80%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "34.84"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "

# Row 15
$ws.Range("D15").Value = "3.639.00"
$ws.Range("E15").Value = "  +0.62%  "

# Row 16
$ws.Range("E16").Value = "  +2.93%  "

# Row 17
$ws.Range("D17").Value = "64.062.80"
$ws.Range("E17").Value = "  +1.34%  "

# Row 18
$ws.Range("D18").Value = "3.128.30"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.79"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.12%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "478.21"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.98%  "

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.43"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.703"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.74%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.59"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.46%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "87.03"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +4.92%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "13.29"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.72"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.00%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "8.22"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.63%  "

# Row 29
$ws.Range("E29").Value = "  +5.70%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.04"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.111"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.72%  "

# Row 32
$ws.Range("E32").Value = "  +0.13%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "26.55"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.46%  "

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.63"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.23%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.08"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -1.08%  "

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.95"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.70%  "

# Row 37
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "52.56"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.36%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0743"
$ws.Range("E38").Value = "  +1.30%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.95"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +2.36%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "434.55"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -3.70%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.0388"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "

# Row 42
$ws.Range("E42").Value = "  +1.35%  "

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "8.16"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.82%  "

# Row 44
$ws.Range("D44").Value = "2.856.07"
$ws.Range("E44").Value = "  +1.23%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.256"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.84%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.19"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.49%  "

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.41"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +1.03%  "

# Row 48
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "25.60"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.112"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.74%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "121.30"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.86%  "

